# Auto-generated: update TPM-derived values per commit 'update scripts wuth new tpm'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 8.459557
$ws.Cells.Item(2, 8).Value = 25.378671
$ws.Cells.Item(2, 9).Value = 0.3030961495696597
$ws.Cells.Item(2, 10).Value = 0.3030961495696597
$ws.Cells.Item(2, 13).Value = 19.21315233333334
$ws.Cells.Item(2, 14).Value = 57.63945700000001
$ws.Cells.Item(2, 15).Value = 0.04451179209991234
$ws.Cells.Item(2, 16).Value = 0.04451179209991233
$ws.Cells.Item(2, 17).Value = 162.5347573135164
$ws.Cells.Item(2, 18).Value = 1462.812815821647
$ws.Cells.Item(2, 19).Value = 0.01349135279592863
$ws.Cells.Item(2, 20).Value = 0.01349135279592862

# Row 3
$ws.Cells.Item(3, 7).Value = 8.459557
$ws.Cells.Item(3, 8).Value = 25.378671
$ws.Cells.Item(3, 9).Value = 0.3030961495696597
$ws.Cells.Item(3, 10).Value = 0.3030961495696597
$ws.Cells.Item(3, 15).Value = 0.2141755495962477
$ws.Cells.Item(3, 16).Value = 0.2141755495962477
$ws.Cells.Item(3, 17).Value = 782.0617713611147
$ws.Cells.Item(3, 18).Value = 7038.555942250032
$ws.Cells.Item(3, 19).Value = 0.06491578441458838
$ws.Cells.Item(3, 20).Value = 0.06491578441458837

# Row 4
$ws.Cells.Item(4, 7).Value = 8.459557
$ws.Cells.Item(4, 8).Value = 25.378671
$ws.Cells.Item(4, 9).Value = 0.3030961495696597
$ws.Cells.Item(4, 10).Value = 0.3030961495696597
$ws.Cells.Item(4, 13).Value = 166.8580016666666
$ws.Cells.Item(4, 14).Value = 500.5740049999999
$ws.Cells.Item(4, 15).Value = 0.3865658561145097
$ws.Cells.Item(4, 16).Value = 0.3865658561145097
$ws.Cells.Item(4, 17).Value = 1411.544776005261
$ws.Cells.Item(4, 18).Value = 12703.90298404735
$ws.Cells.Item(4, 19).Value = 0.117166622543407
$ws.Cells.Item(4, 20).Value = 0.117166622543407

# Row 5
$ws.Cells.Item(5, 7).Value = 8.459557
$ws.Cells.Item(5, 8).Value = 25.378671
$ws.Cells.Item(5, 9).Value = 0.3030961495696597
$ws.Cells.Item(5, 10).Value = 0.3030961495696597
$ws.Cells.Item(5, 13).Value = 41.09915599999999
$ws.Cells.Item(5, 14).Value = 123.297468
$ws.Cells.Item(5, 15).Value = 0.09521587377309249
$ws.Cells.Item(5, 16).Value = 0.09521587377309249
$ws.Cells.Item(5, 17).Value = 347.6806528338919
$ws.Cells.Item(5, 18).Value = 3129.125875505028
$ws.Cells.Item(5, 19).Value = 0.02885956471853508
$ws.Cells.Item(5, 20).Value = 0.02885956471853508

# Row 6
$ws.Cells.Item(6, 7).Value = 8.459557
$ws.Cells.Item(6, 8).Value = 25.378671
$ws.Cells.Item(6, 9).Value = 0.3030961495696597
$ws.Cells.Item(6, 10).Value = 0.3030961495696597
$ws.Cells.Item(6, 13).Value = 112.0244103333333
$ws.Cells.Item(6, 14).Value = 336.073231
$ws.Cells.Item(6, 15).Value = 0.2595309284162377
$ws.Cells.Item(6, 16).Value = 0.2595309284162377
$ws.Cells.Item(6, 17).Value = 947.6768846062224
$ws.Cells.Item(6, 18).Value = 8529.091961456003
$ws.Cells.Item(6, 19).Value = 0.07866282509720064
$ws.Cells.Item(6, 20).Value = 0.07866282509720063

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4601547065605718
$ws.Cells.Item(7, 10).Value = 0.4601547065605718
$ws.Cells.Item(7, 13).Value = 19.21315233333334
$ws.Cells.Item(7, 14).Value = 57.63945700000001
$ws.Cells.Item(7, 15).Value = 0.04451179209991234
$ws.Cells.Item(7, 16).Value = 0.04451179209991233
$ws.Cells.Item(7, 17).Value = 246.7571220013333
$ws.Cells.Item(7, 18).Value = 2220.814098011999
$ws.Cells.Item(7, 19).Value = 0.02048231063222034
$ws.Cells.Item(7, 20).Value = 0.02048231063222034

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4601547065605718
$ws.Cells.Item(8, 10).Value = 0.4601547065605718
$ws.Cells.Item(8, 15).Value = 0.2141755495962477
$ws.Cells.Item(8, 16).Value = 0.2141755495962477
$ws.Cells.Item(8, 19).Value = 0.09855388717691058
$ws.Cells.Item(8, 20).Value = 0.09855388717691055

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4601547065605718
$ws.Cells.Item(9, 10).Value = 0.4601547065605718
$ws.Cells.Item(9, 13).Value = 166.8580016666666
$ws.Cells.Item(9, 14).Value = 500.5740049999999
$ws.Cells.Item(9, 15).Value = 0.3865658561145097
$ws.Cells.Item(9, 16).Value = 0.3865658561145097
$ws.Cells.Item(9, 17).Value = 2142.979952473892
$ws.Cells.Item(9, 18).Value = 19286.81957226503
$ws.Cells.Item(9, 19).Value = 0.1778800980867085
$ws.Cells.Item(9, 20).Value = 0.1778800980867084

# Row 10
$ws.Cells.Item(10, 9).Value = 0.4601547065605718
$ws.Cells.Item(10, 10).Value = 0.4601547065605718
$ws.Cells.Item(10, 13).Value = 41.09915599999999
$ws.Cells.Item(10, 14).Value = 123.297468
$ws.Cells.Item(10, 15).Value = 0.09521587377309249
$ws.Cells.Item(10, 16).Value = 0.09521587377309249
$ws.Cells.Item(10, 17).Value = 527.8420362934972
$ws.Cells.Item(10, 18).Value = 4750.578326641476
$ws.Cells.Item(10, 19).Value = 0.04381403245596582
$ws.Cells.Item(10, 20).Value = 0.04381403245596582

# Row 11
$ws.Cells.Item(11, 9).Value = 0.4601547065605718
$ws.Cells.Item(11, 10).Value = 0.4601547065605718
$ws.Cells.Item(11, 13).Value = 112.0244103333333
$ws.Cells.Item(11, 14).Value = 336.073231
$ws.Cells.Item(11, 15).Value = 0.2595309284162377
$ws.Cells.Item(11, 16).Value = 0.2595309284162377
$ws.Cells.Item(11, 17).Value = 1438.744699889335
$ws.Cells.Item(11, 18).Value = 12948.70229900402
$ws.Cells.Item(11, 19).Value = 0.1194243782087666
$ws.Cells.Item(11, 20).Value = 0.1194243782087666

# Row 12
$ws.Cells.Item(12, 7).Value = 1.955432333333333
$ws.Cells.Item(12, 8).Value = 5.866296999999999
$ws.Cells.Item(12, 9).Value = 0.0700608803720276
$ws.Cells.Item(12, 10).Value = 0.0700608803720276
$ws.Cells.Item(12, 13).Value = 19.21315233333334
$ws.Cells.Item(12, 14).Value = 57.63945700000001
$ws.Cells.Item(12, 15).Value = 0.04451179209991234
$ws.Cells.Item(12, 16).Value = 0.04451179209991233
$ws.Cells.Item(12, 17).Value = 37.57001929785879
$ws.Cells.Item(12, 18).Value = 338.130173680729
$ws.Cells.Item(12, 19).Value = 0.003118535341456521
$ws.Cells.Item(12, 20).Value = 0.003118535341456521

# Row 13
$ws.Cells.Item(13, 7).Value = 1.955432333333333
$ws.Cells.Item(13, 8).Value = 5.866296999999999
$ws.Cells.Item(13, 9).Value = 0.0700608803720276
$ws.Cells.Item(13, 10).Value = 0.0700608803720276
$ws.Cells.Item(13, 15).Value = 0.2141755495962477
$ws.Cells.Item(13, 16).Value = 0.2141755495962477
$ws.Cells.Item(13, 17).Value = 180.7741084294915
$ws.Cells.Item(13, 18).Value = 1626.966975865424
$ws.Cells.Item(13, 19).Value = 0.01500532755887598
$ws.Cells.Item(13, 20).Value = 0.01500532755887597

# Row 14
$ws.Cells.Item(14, 7).Value = 1.955432333333333
$ws.Cells.Item(14, 8).Value = 5.866296999999999
$ws.Cells.Item(14, 9).Value = 0.0700608803720276
$ws.Cells.Item(14, 10).Value = 0.0700608803720276
$ws.Cells.Item(14, 13).Value = 166.8580016666666
$ws.Cells.Item(14, 14).Value = 500.5740049999999
$ws.Cells.Item(14, 15).Value = 0.3865658561145097
$ws.Cells.Item(14, 16).Value = 0.3865658561145097
$ws.Cells.Item(14, 17).Value = 326.2795315343872
$ws.Cells.Item(14, 18).Value = 2936.515783809484
$ws.Cells.Item(14, 19).Value = 0.0270831442011491
$ws.Cells.Item(14, 20).Value = 0.0270831442011491

# Row 15
$ws.Cells.Item(15, 7).Value = 1.955432333333333
$ws.Cells.Item(15, 8).Value = 5.866296999999999
$ws.Cells.Item(15, 9).Value = 0.0700608803720276
$ws.Cells.Item(15, 10).Value = 0.0700608803720276
$ws.Cells.Item(15, 13).Value = 41.09915599999999
$ws.Cells.Item(15, 14).Value = 123.297468
$ws.Cells.Item(15, 15).Value = 0.09521587377309249
$ws.Cells.Item(15, 16).Value = 0.09521587377309249
$ws.Cells.Item(15, 17).Value = 80.36661851511064
$ws.Cells.Item(15, 18).Value = 723.2995666359958
$ws.Cells.Item(15, 19).Value = 0.006670907941934713
$ws.Cells.Item(15, 20).Value = 0.006670907941934713

# Row 16
$ws.Cells.Item(16, 7).Value = 1.955432333333333
$ws.Cells.Item(16, 8).Value = 5.866296999999999
$ws.Cells.Item(16, 9).Value = 0.0700608803720276
$ws.Cells.Item(16, 10).Value = 0.0700608803720276
$ws.Cells.Item(16, 13).Value = 112.0244103333333
$ws.Cells.Item(16, 14).Value = 336.073231
$ws.Cells.Item(16, 15).Value = 0.2595309284162377
$ws.Cells.Item(16, 16).Value = 0.2595309284162377
$ws.Cells.Item(16, 17).Value = 219.0561540884008
$ws.Cells.Item(16, 18).Value = 1971.505386795607
$ws.Cells.Item(16, 19).Value = 0.01818296532861129
$ws.Cells.Item(16, 20).Value = 0.01818296532861129

# Row 17
$ws.Cells.Item(17, 7).Value = 2.929608
$ws.Cells.Item(17, 8).Value = 8.788824000000002
$ws.Cells.Item(17, 9).Value = 0.1049644685352285
$ws.Cells.Item(17, 10).Value = 0.1049644685352285
$ws.Cells.Item(17, 13).Value = 19.21315233333334
$ws.Cells.Item(17, 14).Value = 57.63945700000001
$ws.Cells.Item(17, 15).Value = 0.04451179209991234
$ws.Cells.Item(17, 16).Value = 0.04451179209991233
$ws.Cells.Item(17, 17).Value = 56.28700478095202
$ws.Cells.Item(17, 18).Value = 506.5830430285682
$ws.Cells.Item(17, 19).Value = 0.00467215660131788
$ws.Cells.Item(17, 20).Value = 0.004672156601317879

# Row 18
$ws.Cells.Item(18, 7).Value = 2.929608
$ws.Cells.Item(18, 8).Value = 8.788824000000002
$ws.Cells.Item(18, 9).Value = 0.1049644685352285
$ws.Cells.Item(18, 10).Value = 0.1049644685352285
$ws.Cells.Item(18, 15).Value = 0.2141755495962477
$ws.Cells.Item(18, 16).Value = 0.2141755495962477
$ws.Cells.Item(18, 17).Value = 270.833853578112
$ws.Cells.Item(18, 18).Value = 2437.504682203009
$ws.Cells.Item(18, 19).Value = 0.02248082273661061
$ws.Cells.Item(18, 20).Value = 0.02248082273661061

# Row 19
$ws.Cells.Item(19, 7).Value = 2.929608
$ws.Cells.Item(19, 8).Value = 8.788824000000002
$ws.Cells.Item(19, 9).Value = 0.1049644685352285
$ws.Cells.Item(19, 10).Value = 0.1049644685352285
$ws.Cells.Item(19, 13).Value = 166.8580016666666
$ws.Cells.Item(19, 14).Value = 500.5740049999999
$ws.Cells.Item(19, 15).Value = 0.3865658561145097
$ws.Cells.Item(19, 16).Value = 0.3865658561145097
$ws.Cells.Item(19, 17).Value = 488.82853654668
$ws.Cells.Item(19, 18).Value = 4399.456828920121
$ws.Cells.Item(19, 19).Value = 0.04057567964092511
$ws.Cells.Item(19, 20).Value = 0.04057567964092511

# Row 20
$ws.Cells.Item(20, 7).Value = 2.929608
$ws.Cells.Item(20, 8).Value = 8.788824000000002
$ws.Cells.Item(20, 9).Value = 0.1049644685352285
$ws.Cells.Item(20, 10).Value = 0.1049644685352285
$ws.Cells.Item(20, 13).Value = 41.09915599999999
$ws.Cells.Item(20, 14).Value = 123.297468
$ws.Cells.Item(20, 15).Value = 0.09521587377309249
$ws.Cells.Item(20, 16).Value = 0.09521587377309249
$ws.Cells.Item(20, 17).Value = 120.404416210848
$ws.Cells.Item(20, 18).Value = 1083.639745897632
$ws.Cells.Item(20, 19).Value = 0.009994283586710053
$ws.Cells.Item(20, 20).Value = 0.009994283586710053

# Row 21
$ws.Cells.Item(21, 7).Value = 2.929608
$ws.Cells.Item(21, 8).Value = 8.788824000000002
$ws.Cells.Item(21, 9).Value = 0.1049644685352285
$ws.Cells.Item(21, 10).Value = 0.1049644685352285
$ws.Cells.Item(21, 13).Value = 112.0244103333333
$ws.Cells.Item(21, 14).Value = 336.073231
$ws.Cells.Item(21, 15).Value = 0.2595309284162377
$ws.Cells.Item(21, 16).Value = 0.2595309284162377
$ws.Cells.Item(21, 17).Value = 328.1876087078161
$ws.Cells.Item(21, 18).Value = 2953.688478370345
$ws.Cells.Item(21, 19).Value = 0.02724152596966482
$ws.Cells.Item(21, 20).Value = 0.02724152596966482

# Row 22
$ws.Cells.Item(22, 7).Value = 1.722740333333333
$ws.Cells.Item(22, 8).Value = 5.168221
$ws.Cells.Item(22, 9).Value = 0.06172379496251228
$ws.Cells.Item(22, 10).Value = 0.06172379496251227
$ws.Cells.Item(22, 13).Value = 19.21315233333334
$ws.Cells.Item(22, 14).Value = 57.63945700000001
$ws.Cells.Item(22, 15).Value = 0.04451179209991234
$ws.Cells.Item(22, 16).Value = 0.04451179209991233
$ws.Cells.Item(22, 17).Value = 33.09927245511079
$ws.Cells.Item(22, 18).Value = 297.893452095997
$ws.Cells.Item(22, 19).Value = 0.002747436728988963
$ws.Cells.Item(22, 20).Value = 0.002747436728988962

# Row 23
$ws.Cells.Item(23, 7).Value = 1.722740333333333
$ws.Cells.Item(23, 8).Value = 5.168221
$ws.Cells.Item(23, 9).Value = 0.06172379496251228
$ws.Cells.Item(23, 10).Value = 0.06172379496251227
$ws.Cells.Item(23, 15).Value = 0.2141755495962477
$ws.Cells.Item(23, 16).Value = 0.2141755495962477
$ws.Cells.Item(23, 17).Value = 159.2624007004036
$ws.Cells.Item(23, 18).Value = 1433.361606303632
$ws.Cells.Item(23, 19).Value = 0.01321972770926217
$ws.Cells.Item(23, 20).Value = 0.01321972770926217

# Row 24
$ws.Cells.Item(24, 7).Value = 1.722740333333333
$ws.Cells.Item(24, 8).Value = 5.168221
$ws.Cells.Item(24, 9).Value = 0.06172379496251228
$ws.Cells.Item(24, 10).Value = 0.06172379496251227
$ws.Cells.Item(24, 13).Value = 166.8580016666666
$ws.Cells.Item(24, 14).Value = 500.5740049999999
$ws.Cells.Item(24, 15).Value = 0.3865658561145097
$ws.Cells.Item(24, 16).Value = 0.3865658561145097
$ws.Cells.Item(24, 17).Value = 287.4530094105672
$ws.Cells.Item(24, 18).Value = 2587.077084695105
$ws.Cells.Item(24, 19).Value = 0.02386031164232002
$ws.Cells.Item(24, 20).Value = 0.02386031164232002

# Row 25
$ws.Cells.Item(25, 7).Value = 1.722740333333333
$ws.Cells.Item(25, 8).Value = 5.168221
$ws.Cells.Item(25, 9).Value = 0.06172379496251228
$ws.Cells.Item(25, 10).Value = 0.06172379496251227
$ws.Cells.Item(25, 13).Value = 41.09915599999999
$ws.Cells.Item(25, 14).Value = 123.297468
$ws.Cells.Item(25, 15).Value = 0.09521587377309249
$ws.Cells.Item(25, 16).Value = 0.09521587377309249
$ws.Cells.Item(25, 17).Value = 70.80317370715865
$ws.Cells.Item(25, 18).Value = 637.2285633644279
$ws.Cells.Item(25, 19).Value = 0.005877085069946811
$ws.Cells.Item(25, 20).Value = 0.00587708506994681

# Row 26
$ws.Cells.Item(26, 7).Value = 1.722740333333333
$ws.Cells.Item(26, 8).Value = 5.168221
$ws.Cells.Item(26, 9).Value = 0.06172379496251228
$ws.Cells.Item(26, 10).Value = 0.06172379496251227
$ws.Cells.Item(26, 13).Value = 112.0244103333333
$ws.Cells.Item(26, 14).Value = 336.073231
$ws.Cells.Item(26, 15).Value = 0.2595309284162377
$ws.Cells.Item(26, 16).Value = 0.2595309284162377
$ws.Cells.Item(26, 17).Value = 192.9889699991168
$ws.Cells.Item(26, 18).Value = 1736.900729992051
$ws.Cells.Item(26, 19).Value = 0.01601923381199431
$ws.Cells.Item(26, 20).Value = 0.01601923381199431
